# Updates the cryptos list figures (price / 1h volume change) to the latest
# scraped values, as produced by the scheduled GitHub Actions job.
#
# Numeric-looking strings in the "Price" column must remain plain text
# (the column holds free-form, locale-formatted numbers such as
# "26.605.95" or values with significant trailing zeros like "2.40"),
# so cells whose new value could otherwise be auto-parsed by Excel as a
# number are written with a leading apostrophe to force text storage and
# preserve the exact formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '26.605.95'
$ws.Range("E2").Value = '  +0.75%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '1.640.89'
$ws.Range("E3").Value = '  +1.07%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  +0.11%  '

# Row 5 - BNB
$ws.Range("D5").Value = "'214.75"
$ws.Range("E5").Value = '  +0.99%  '

# Row 6 - XRP
$ws.Range("E6").Value = '  +1.55%  '

# Row 7 - USDC
$ws.Range("E7").Value = '  -0.06%  '

# Row 8 - Cardano
$ws.Range("E8").Value = '  +0.87%  '

# Row 9 - Dogecoin
$ws.Range("E9").Value = '  +0.54%  '

# Row 10 - Solana
$ws.Range("D10").Value = "'19.09"
$ws.Range("E10").Value = '  +0.83%  '

# Row 11 - TRON
$ws.Range("D11").Value = "'0.0841"
$ws.Range("E11").Value = '  -0.12%  '

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = '1.869.16'
$ws.Range("E12").Value = '  +1.08%  '

# Row 13 and 14 swapped places (WrappedEther now ranks above Polkadot)
# Row 13 - now WrappedEther
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.656.85'
$ws.Range("E13").Value = '  +2.75%  '

# Row 14 - now Polkadot
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = "'4.18"
$ws.Range("E14").Value = '  +1.48%  '

# Row 15 - Polygon
$ws.Range("E15").Value = '  +1.44%  '

# Row 16 - Litecoin
$ws.Range("D16").Value = "'64.82"
$ws.Range("E16").Value = '  +1.02%  '

# Row 17 - WrappedBTC
$ws.Range("D17").Value = '26.611.97'
$ws.Range("E17").Value = '  +0.70%  '

# Row 18 - ShibaInu
$ws.Range("E18").Value = '  +0.47%  '

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'215.36"
$ws.Range("E19").Value = '  +0.19%  '

# Row 20 - Dai
$ws.Range("E20").Value = '  -0.14%  '

# Row 21 - Uniswap
$ws.Range("E21").Value = '  +0.91%  '

# Row 22 - Chainlink
$ws.Range("E22").Value = '  +0.45%  '

# Row 23 - Avalanche
$ws.Range("E23").Value = '  +1.73%  '

# Row 24 - Toncoin
$ws.Range("D24").Value = "'2.21"
$ws.Range("E24").Value = '  +12.38%  '

# Row 25 - Monero
$ws.Range("E25").Value = '  -1.80%  '

# Row 26 - BinanceUSD
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = '  +0.23%  '

# Row 27 - Stellar
$ws.Range("E27").Value = '  +0.14%  '

# Row 28 - Cosmos
$ws.Range("E28").Value = '  +4.04%  '

# Row 29 - EthereumClassic
$ws.Range("E29").Value = '  +0.81%  '

# Row 30 - Hedera
$ws.Range("E30").Value = '  +0.99%  '

# Row 31 - PancakeSwap
$ws.Range("E31").Value = '  +1.66%  '

# Row 32 - Filecoin
$ws.Range("E32").Value = '  +1.31%  '

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "'2.99"
$ws.Range("E33").Value = '  +1.61%  '

# Row 34 - Maker
$ws.Range("D34").Value = '1.276.77'
$ws.Range("E34").Value = '  +5.06%  '

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = '  +3.02%  '

# Row 36 - HuobiToken
$ws.Range("D36").Value = "'2.40"
$ws.Range("E36").Value = '  +1.13%  '

# Row 37 - VeChain
$ws.Range("E37").Value = '  +2.69%  '

# Row 38 - ImmutableX
$ws.Range("E38").Value = '  +6.07%  '

# Row 39 - ARBITRUM
$ws.Range("D39").Value = "'0.823"
$ws.Range("E39").Value = '  +3.49%  '

# Row 40 - PaxDollar
$ws.Range("E40").Value = '  -0.13%  '

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "'0.809"
$ws.Range("E41").Value = '  +2.16%  '

# Row 43 - FraxShare
$ws.Range("E43").Value = '  +0.92%  '

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = '1.779.40'
$ws.Range("E44").Value = '  +1.09%  '

# Row 45 - Quant
$ws.Range("D45").Value = "'91.50"
$ws.Range("E45").Value = '  -1.31%  '

# Row 46 - Aave
$ws.Range("D46").Value = "'59.11"
$ws.Range("E46").Value = '  +8.04%  '

# Row 47 - RenderToken
$ws.Range("E47").Value = '  +0.76%  '

# Row 48 - Cronos
$ws.Range("E48").Value = '  +0.99%  '

# Row 49 - EnergySwap
$ws.Range("E49").Value = '  +1.62%  '

# Row 50 - Algorand
$ws.Range("E50").Value = '  +1.21%  '

# Row 51 - Mantle
$ws.Range("E51").Value = '  -0.32%  '
